$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1748284055722/p2gySLFM9k-camera-video-segment-driver-1748284058222.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDCWD5SDZ%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T150021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIGoXMSY%2FyPtoC0htMFLZeLI6EHbQS2Gu2x5atnbcS333AiEA4yVEB9DAUa8fOVai61HoAvC1k67ndil%2F9oEhKyG3Aloq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDBJu%2BHdjOTO7NBcGiyq6A21PcI1%2BoSpzeEaSzEaeLnsBUo6R3Ysd8qmLcKqn0ZeDPVl6Qp%2FWrTxCrEP1A35e9Yc64AE6F3EkYclxwS%2BO%2Bevun7Z2nKCF%2FqFYuqyhcy%2FADTj2tGDiPA%2Bd18sWydiiCCBznzIARbAsRnVHtbZE3JhFtyFqbUpVaCfnvYTYDdvWkRUlXCwiYBboHnc0XDdYPTSfo3rtpVWXvDa8OpLOGv%2BKnu7wPkqAvLj4%2B8DuAcWJzeggoV2Zh4jVnHcawkKcP72ZMmJMHvrFyLMQEmltImOhwJY2olqcpcxmB5dQcCh3eZQP1Vbt4FTvydsjhVOVEnsi2Mn25FVTldpjSZPjbqPqwt2%2B3T0wA8DHa54zemOA7OcOlnoMGTVSrfhSl5O4oLKzilz5EPUcRMjyLvWMy502zByXVujkZq7MhlMsPfxD1z%2FGy07vL5HxIHxEsfFx3xRlA18JjuTGGEufLT5K9P5GxQQR%2Fc4CQ%2F5fSG0qqhVgtOJH%2BtfWagsv8Wh16XCYIxqW5IOL7GcdUW%2Fl6fF8SGSt1g40NTxFd5JwdXG7MBAwqH46GQL1FWrdoYhIRwEpC8pCjvYOY%2BTcOtIw%2F4DXwQY6pQG3VQom66%2Bgm8nszWyt6AeVH8n319cmtrIllETf06KoLyqtYau5yMLWYUO25A8e75vU0PUYYDLj%2BaOQckoLdDPSvPA7vey1R4COxgzTX18b03QzGUgXx6Xe%2Bsn3dKaVYGFS%2FBgPhPopKI3jPVEOvVB9f%2Fvy1nK32W7%2BUAB05QDbiXG4iDR8MGzQ0A%2FnJJmi3YvFZUbSlzT3UaUK5L0l5SimABeoElE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2023%3A00%3A21%20GMT&X-Amz-Signature=8322147f9f7ffa73d21529fd9217e6eadfc2105db4d520e2be22e01e1c727a9a'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1748284055722/JuEVYDqJVY-camera-video-segment-1748284058222.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDCWD5SDZ%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T150021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIGoXMSY%2FyPtoC0htMFLZeLI6EHbQS2Gu2x5atnbcS333AiEA4yVEB9DAUa8fOVai61HoAvC1k67ndil%2F9oEhKyG3Aloq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDBJu%2BHdjOTO7NBcGiyq6A21PcI1%2BoSpzeEaSzEaeLnsBUo6R3Ysd8qmLcKqn0ZeDPVl6Qp%2FWrTxCrEP1A35e9Yc64AE6F3EkYclxwS%2BO%2Bevun7Z2nKCF%2FqFYuqyhcy%2FADTj2tGDiPA%2Bd18sWydiiCCBznzIARbAsRnVHtbZE3JhFtyFqbUpVaCfnvYTYDdvWkRUlXCwiYBboHnc0XDdYPTSfo3rtpVWXvDa8OpLOGv%2BKnu7wPkqAvLj4%2B8DuAcWJzeggoV2Zh4jVnHcawkKcP72ZMmJMHvrFyLMQEmltImOhwJY2olqcpcxmB5dQcCh3eZQP1Vbt4FTvydsjhVOVEnsi2Mn25FVTldpjSZPjbqPqwt2%2B3T0wA8DHa54zemOA7OcOlnoMGTVSrfhSl5O4oLKzilz5EPUcRMjyLvWMy502zByXVujkZq7MhlMsPfxD1z%2FGy07vL5HxIHxEsfFx3xRlA18JjuTGGEufLT5K9P5GxQQR%2Fc4CQ%2F5fSG0qqhVgtOJH%2BtfWagsv8Wh16XCYIxqW5IOL7GcdUW%2Fl6fF8SGSt1g40NTxFd5JwdXG7MBAwqH46GQL1FWrdoYhIRwEpC8pCjvYOY%2BTcOtIw%2F4DXwQY6pQG3VQom66%2Bgm8nszWyt6AeVH8n319cmtrIllETf06KoLyqtYau5yMLWYUO25A8e75vU0PUYYDLj%2BaOQckoLdDPSvPA7vey1R4COxgzTX18b03QzGUgXx6Xe%2Bsn3dKaVYGFS%2FBgPhPopKI3jPVEOvVB9f%2Fvy1nK32W7%2BUAB05QDbiXG4iDR8MGzQ0A%2FnJJmi3YvFZUbSlzT3UaUK5L0l5SimABeoElE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2023%3A00%3A21%20GMT&X-Amz-Signature=97ef5c719c48f75d5982984aec7220ab1b7938d90587c448027da2fb2e88d9bf'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1748283979463/3bY9s5rKOY-camera-video-segment-driver-1748283981963.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDCWD5SDZ%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T150021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIGoXMSY%2FyPtoC0htMFLZeLI6EHbQS2Gu2x5atnbcS333AiEA4yVEB9DAUa8fOVai61HoAvC1k67ndil%2F9oEhKyG3Aloq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDBJu%2BHdjOTO7NBcGiyq6A21PcI1%2BoSpzeEaSzEaeLnsBUo6R3Ysd8qmLcKqn0ZeDPVl6Qp%2FWrTxCrEP1A35e9Yc64AE6F3EkYclxwS%2BO%2Bevun7Z2nKCF%2FqFYuqyhcy%2FADTj2tGDiPA%2Bd18sWydiiCCBznzIARbAsRnVHtbZE3JhFtyFqbUpVaCfnvYTYDdvWkRUlXCwiYBboHnc0XDdYPTSfo3rtpVWXvDa8OpLOGv%2BKnu7wPkqAvLj4%2B8DuAcWJzeggoV2Zh4jVnHcawkKcP72ZMmJMHvrFyLMQEmltImOhwJY2olqcpcxmB5dQcCh3eZQP1Vbt4FTvydsjhVOVEnsi2Mn25FVTldpjSZPjbqPqwt2%2B3T0wA8DHa54zemOA7OcOlnoMGTVSrfhSl5O4oLKzilz5EPUcRMjyLvWMy502zByXVujkZq7MhlMsPfxD1z%2FGy07vL5HxIHxEsfFx3xRlA18JjuTGGEufLT5K9P5GxQQR%2Fc4CQ%2F5fSG0qqhVgtOJH%2BtfWagsv8Wh16XCYIxqW5IOL7GcdUW%2Fl6fF8SGSt1g40NTxFd5JwdXG7MBAwqH46GQL1FWrdoYhIRwEpC8pCjvYOY%2BTcOtIw%2F4DXwQY6pQG3VQom66%2Bgm8nszWyt6AeVH8n319cmtrIllETf06KoLyqtYau5yMLWYUO25A8e75vU0PUYYDLj%2BaOQckoLdDPSvPA7vey1R4COxgzTX18b03QzGUgXx6Xe%2Bsn3dKaVYGFS%2FBgPhPopKI3jPVEOvVB9f%2Fvy1nK32W7%2BUAB05QDbiXG4iDR8MGzQ0A%2FnJJmi3YvFZUbSlzT3UaUK5L0l5SimABeoElE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2023%3A00%3A21%20GMT&X-Amz-Signature=ce9634012251652071220fcab5bd52061b8896a78d5ae4a67723685fba0a10b2'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1748283979463/qChm84icV1-camera-video-segment-1748283981963.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDCWD5SDZ%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T150021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIGoXMSY%2FyPtoC0htMFLZeLI6EHbQS2Gu2x5atnbcS333AiEA4yVEB9DAUa8fOVai61HoAvC1k67ndil%2F9oEhKyG3Aloq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDBJu%2BHdjOTO7NBcGiyq6A21PcI1%2BoSpzeEaSzEaeLnsBUo6R3Ysd8qmLcKqn0ZeDPVl6Qp%2FWrTxCrEP1A35e9Yc64AE6F3EkYclxwS%2BO%2Bevun7Z2nKCF%2FqFYuqyhcy%2FADTj2tGDiPA%2Bd18sWydiiCCBznzIARbAsRnVHtbZE3JhFtyFqbUpVaCfnvYTYDdvWkRUlXCwiYBboHnc0XDdYPTSfo3rtpVWXvDa8OpLOGv%2BKnu7wPkqAvLj4%2B8DuAcWJzeggoV2Zh4jVnHcawkKcP72ZMmJMHvrFyLMQEmltImOhwJY2olqcpcxmB5dQcCh3eZQP1Vbt4FTvydsjhVOVEnsi2Mn25FVTldpjSZPjbqPqwt2%2B3T0wA8DHa54zemOA7OcOlnoMGTVSrfhSl5O4oLKzilz5EPUcRMjyLvWMy502zByXVujkZq7MhlMsPfxD1z%2FGy07vL5HxIHxEsfFx3xRlA18JjuTGGEufLT5K9P5GxQQR%2Fc4CQ%2F5fSG0qqhVgtOJH%2BtfWagsv8Wh16XCYIxqW5IOL7GcdUW%2Fl6fF8SGSt1g40NTxFd5JwdXG7MBAwqH46GQL1FWrdoYhIRwEpC8pCjvYOY%2BTcOtIw%2F4DXwQY6pQG3VQom66%2Bgm8nszWyt6AeVH8n319cmtrIllETf06KoLyqtYau5yMLWYUO25A8e75vU0PUYYDLj%2BaOQckoLdDPSvPA7vey1R4COxgzTX18b03QzGUgXx6Xe%2Bsn3dKaVYGFS%2FBgPhPopKI3jPVEOvVB9f%2Fvy1nK32W7%2BUAB05QDbiXG4iDR8MGzQ0A%2FnJJmi3YvFZUbSlzT3UaUK5L0l5SimABeoElE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2023%3A00%3A21%20GMT&X-Amz-Signature=0bb42030c591ef550ec3e3e5784406b8606d381a46028330eac83ec68d303b39'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1748276157617/6zxWnZYxXR-camera-video-segment-driver-1748276160117.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDCWD5SDZ%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T150021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIGoXMSY%2FyPtoC0htMFLZeLI6EHbQS2Gu2x5atnbcS333AiEA4yVEB9DAUa8fOVai61HoAvC1k67ndil%2F9oEhKyG3Aloq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDBJu%2BHdjOTO7NBcGiyq6A21PcI1%2BoSpzeEaSzEaeLnsBUo6R3Ysd8qmLcKqn0ZeDPVl6Qp%2FWrTxCrEP1A35e9Yc64AE6F3EkYclxwS%2BO%2Bevun7Z2nKCF%2FqFYuqyhcy%2FADTj2tGDiPA%2Bd18sWydiiCCBznzIARbAsRnVHtbZE3JhFtyFqbUpVaCfnvYTYDdvWkRUlXCwiYBboHnc0XDdYPTSfo3rtpVWXvDa8OpLOGv%2BKnu7wPkqAvLj4%2B8DuAcWJzeggoV2Zh4jVnHcawkKcP72ZMmJMHvrFyLMQEmltImOhwJY2olqcpcxmB5dQcCh3eZQP1Vbt4FTvydsjhVOVEnsi2Mn25FVTldpjSZPjbqPqwt2%2B3T0wA8DHa54zemOA7OcOlnoMGTVSrfhSl5O4oLKzilz5EPUcRMjyLvWMy502zByXVujkZq7MhlMsPfxD1z%2FGy07vL5HxIHxEsfFx3xRlA18JjuTGGEufLT5K9P5GxQQR%2Fc4CQ%2F5fSG0qqhVgtOJH%2BtfWagsv8Wh16XCYIxqW5IOL7GcdUW%2Fl6fF8SGSt1g40NTxFd5JwdXG7MBAwqH46GQL1FWrdoYhIRwEpC8pCjvYOY%2BTcOtIw%2F4DXwQY6pQG3VQom66%2Bgm8nszWyt6AeVH8n319cmtrIllETf06KoLyqtYau5yMLWYUO25A8e75vU0PUYYDLj%2BaOQckoLdDPSvPA7vey1R4COxgzTX18b03QzGUgXx6Xe%2Bsn3dKaVYGFS%2FBgPhPopKI3jPVEOvVB9f%2Fvy1nK32W7%2BUAB05QDbiXG4iDR8MGzQ0A%2FnJJmi3YvFZUbSlzT3UaUK5L0l5SimABeoElE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2023%3A00%3A21%20GMT&X-Amz-Signature=0b8c1cd5d044da5297aad11712016f9651ba578bfb532d241f8515a201455ee2'
$ws.Range("L4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1748276157617/AqojMFPYnw-camera-video-segment-1748276160117.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSDCWD5SDZ%2F20250527%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250527T150021Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJb%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIGoXMSY%2FyPtoC0htMFLZeLI6EHbQS2Gu2x5atnbcS333AiEA4yVEB9DAUa8fOVai61HoAvC1k67ndil%2F9oEhKyG3Aloq3QMIXxAEGgw3ODEyMDQ5NDIyNDQiDBJu%2BHdjOTO7NBcGiyq6A21PcI1%2BoSpzeEaSzEaeLnsBUo6R3Ysd8qmLcKqn0ZeDPVl6Qp%2FWrTxCrEP1A35e9Yc64AE6F3EkYclxwS%2BO%2Bevun7Z2nKCF%2FqFYuqyhcy%2FADTj2tGDiPA%2Bd18sWydiiCCBznzIARbAsRnVHtbZE3JhFtyFqbUpVaCfnvYTYDdvWkRUlXCwiYBboHnc0XDdYPTSfo3rtpVWXvDa8OpLOGv%2BKnu7wPkqAvLj4%2B8DuAcWJzeggoV2Zh4jVnHcawkKcP72ZMmJMHvrFyLMQEmltImOhwJY2olqcpcxmB5dQcCh3eZQP1Vbt4FTvydsjhVOVEnsi2Mn25FVTldpjSZPjbqPqwt2%2B3T0wA8DHa54zemOA7OcOlnoMGTVSrfhSl5O4oLKzilz5EPUcRMjyLvWMy502zByXVujkZq7MhlMsPfxD1z%2FGy07vL5HxIHxEsfFx3xRlA18JjuTGGEufLT5K9P5GxQQR%2Fc4CQ%2F5fSG0qqhVgtOJH%2BtfWagsv8Wh16XCYIxqW5IOL7GcdUW%2Fl6fF8SGSt1g40NTxFd5JwdXG7MBAwqH46GQL1FWrdoYhIRwEpC8pCjvYOY%2BTcOtIw%2F4DXwQY6pQG3VQom66%2Bgm8nszWyt6AeVH8n319cmtrIllETf06KoLyqtYau5yMLWYUO25A8e75vU0PUYYDLj%2BaOQckoLdDPSvPA7vey1R4COxgzTX18b03QzGUgXx6Xe%2Bsn3dKaVYGFS%2FBgPhPopKI3jPVEOvVB9f%2Fvy1nK32W7%2BUAB05QDbiXG4iDR8MGzQ0A%2FnJJmi3YvFZUbSlzT3UaUK5L0l5SimABeoElE%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2027%20May%202025%2023%3A00%3A21%20GMT&X-Amz-Signature=d5e86c98b5d9a9f197a5f9818ce989e3cda36e9ff98e60981e1d2d276dcfeffb'
